$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Mariam Al Maktoum"
$wsSummary.Range("B4").Value = 5027.59
$wsSummary.Range("B6").Value = 595651
$wsSummary.Range("B7").Value = 375418
$wsSummary.Range("B8").Value = 220233
$wsSummary.Range("B9").Value = 1.59

# --- Assets sheet ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("B2").Value = "Luxury Car"
$wsAssets.Range("C2").Value = 590676
$wsAssets.Range("C3").Value = 4975
$wsAssets.Range("C4").Value = 595651

# --- Liabilities sheet ---
$wsLiabilities = $wb.Worksheets.Item("Liabilities")
# Remove the "Personal Loans" row entirely (row 3); rows below shift up.
$wsLiabilities.Rows(3).Delete()

# Update remaining values after the shift.
$wsLiabilities.Range("C2").Value = 354406
$wsLiabilities.Range("D2").Value = 4922
$wsLiabilities.Range("A3").Value = "Credit Cards"
$wsLiabilities.Range("B3").Value = "Credit Card Balance"
$wsLiabilities.Range("C3").Value = 21012
$wsLiabilities.Range("D3").Value = 1051
$wsLiabilities.Range("E3").Value = 1
$wsLiabilities.Range("C4").Value = 375418
